$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates reflecting the latest cryptocurrency market data refresh.
# Numeric-looking text values are prefixed with a leading apostrophe to force
# Excel to keep them as text (preserving formats like trailing zeros), then the
# cell style is reset to "Normal" to avoid leaving a stray text-number-format.

$ws.Range("D2").Value = '26.379.87'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '1.594.31'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  -0.72%  '
$ws.Range("D5").Value = "'210.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").Value = "'0.505"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("D8").Value = "'0.0611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.82%  '
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").Value = "'19.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").Value = "'0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.77%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.570.94'
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '26.379.65'
$ws.Range("E17").Value = '  -0.90%  '
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").Value = "'7.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.14%  '
$ws.Range("D20").Value = "'211.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  -4.12%  '
$ws.Range("D24").Value = "'8.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = "'144.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("D27").Value = "'7.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").Value = "'1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("D34").Value = '1.314.16'
$ws.Range("E34").Value = '  +2.94%  '
$ws.Range("D35").Value = "'0.617"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.17%  '
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("E37").Value = '  -0.72%  '
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("D39").Value = "'1.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -13.13%  '
$ws.Range("D40").Value = "'0.812"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.00%  '
$ws.Range("E41").Value = '  -0.67%  '
$ws.Range("E42").Value = '  +3.86%  '
$ws.Range("D43").Value = "'2.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = "'0.764"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.63%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = "'62.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = '1.727.02'
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").Value = "'88.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = "'1.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.57%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0102'
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("D50").Value = "'0.0984"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.09%  '
$ws.Range("D51").Value = "'0.0505"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.39%  '
